$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Test questions as nodes" - clean up / rename a few Topic (my opinion) entries,
# replacing curly apostrophes with straight ones and renaming the Nobel-prize
# topic label.

$ws.Range("C4").Value = "Firm's behavior"
$ws.Range("C5").Value = "Consumer's choice"

$ws.Range("C18").Value = "Nobel research (connectivity)"
$ws.Range("C18").Font.Name = "Arial"
$ws.Range("C18").Font.ThemeColor = 1

# Move the active selection to C4 (matches the saved workbook's cursor position)
$ws.Range("C4").Select() | Out-Null
